$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update day labels (capitalized) and temperature values for the first 5 rows
$ws.Range("A1").Value = "Lunes"
$ws.Range("B1").Value = 26.2

$ws.Range("A2").Value = "Martes"
$ws.Range("B2").Value = 21.4

$ws.Range("A3").Value = "Miercoles"
$ws.Range("B3").Value = 20.6

$ws.Range("A4").Value = "Jueves"
$ws.Range("B4").Value = 23.4

$ws.Range("A5").Value = "Viernes"
$ws.Range("B5").Value = 19.8

# Remove the now-unused rows 6 (sabado) and 7 (domingo)
$ws.Range("A6:B7").Delete()

# Update selection to match new last data row
$ws.Range("B5").Select()
